# Week 38 read-me edits: apply changes from the bottom of the
# document upward so earlier paragraph indices stay valid while the
# later replacements change the total paragraph count.

$d = $word.ActiveDocument

# --- Paragraphs 8-10 ("They say you can't compare..." / "Many questions
# remain..." / "I hope you enjoy...") collapse into two new paragraphs
# ("It's been said..." and "I have questions...").
$start = $d.Paragraphs.Item(8).Range.Start
$end = $d.Paragraphs.Item(10).Range.End
$rng = $d.Range($start, $end)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
<w:r><w:t>It</w:t></w:r>
<w:r><w:t>’</w:t></w:r>
<w:r><w:t xml:space="preserve">s been </w:t></w:r>
<w:r><w:t xml:space="preserve">said that athletes </w:t></w:r>
<w:r><w:t xml:space="preserve">from different decades </w:t></w:r>
<w:r><w:t>can't be fairly compared</w:t></w:r>
<w:r><w:t xml:space="preserve">, however </w:t></w:r>
<w:r><w:t xml:space="preserve">this dashboard allows for it. The timeline on the left </w:t></w:r>
<w:r><w:t>plots rank by</w:t></w:r>
<w:r><w:t xml:space="preserve"> calendar years, while the timeline on the right </w:t></w:r>
<w:r><w:t xml:space="preserve">plots rank by </w:t></w:r>
<w:r><w:t>career years.</w:t></w:r>
<w:r><w:t xml:space="preserve"> The dataset starts career years at </w:t></w:r>
<w:r><w:t>0</w:t></w:r>
<w:r><w:t>, I changed to start at 1.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">I have </w:t></w:r>
<w:r><w:t xml:space="preserve">questions about what the ranks truly </w:t></w:r>
<w:r><w:t>represent</w:t></w:r>
<w:r><w:t xml:space="preserve"> and</w:t></w:r>
<w:r><w:t xml:space="preserve"> would love to see a data set like this on</w:t></w:r>
<w:r><w:t>e</w:t></w:r>
<w:r><w:t xml:space="preserve"> covering </w:t></w:r>
<w:r><w:t>a far greater number of years.</w:t></w:r>
</w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng.InsertXML($xml)

# --- Paragraph 7: "Usually, I turn off all grid lines..." -> "Normally, I turn off all grid lines..."
$rng = $d.Paragraphs.Item(7).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
<w:r><w:t xml:space="preserve">Normally, I turn off all grid lines, but </w:t></w:r>
<w:r><w:t>with</w:t></w:r>
<w:r><w:t xml:space="preserve"> the y-axis logarithmic (</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>log_y</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>=True), I decided to keep the horizontal grid lines for better</w:t></w:r>
<w:r><w:t xml:space="preserve"> clarity</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
</w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng.InsertXML($xml)

# --- Paragraph 6: "Dash Mantine ChipGroup was used..." -> "For the first time I used the Dash Mantine ChipGroup..."
$rng = $d.Paragraphs.Item(6).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
<w:r><w:t xml:space="preserve">For the first time </w:t></w:r>
<w:r><w:t xml:space="preserve">I used the Dash Mantine </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ChipGroup</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> for league selection (allowing single or multiple selections). </w:t></w:r>
<w:r><w:t>I</w:t></w:r>
<w:r><w:t>t was a great learning experience</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
</w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng.InsertXML($xml)

# --- Paragraph 5: "Athletes within each league were sorted..." -> condensed single-run version
$rng = $d.Paragraphs.Item(5).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
<w:r><w:t xml:space="preserve">Athletes within each league were first sorted by median rank (to minimize the influence of outliers) across all their years, then by mean rank. The top 5 athletes from each league were selected based on these criteria. </w:t></w:r>
</w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng.InsertXML($xml)

# --- Paragraph 4: "Athletes with less than 5 years..." -> "Athletes with fewer than 5 years... were excluded."
$rng = $d.Paragraphs.Item(4).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
<w:r><w:t>Athletes with fewer than 5 years of data were excluded</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
</w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng.InsertXML($xml)

# --- Paragraph 3: "He" + "re are the details:" -> "He" + "re is a dashboard to look at athlete rankings over time. A few points:"
$rng = $d.Paragraphs.Item(3).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>He</w:t></w:r>
<w:r><w:t xml:space="preserve">re </w:t></w:r>
<w:r><w:t xml:space="preserve">is a dashboard </w:t></w:r>
<w:r><w:t>to look at athlete rankings over time. A few points</w:t></w:r>
<w:r><w:t>:</w:t></w:r>
</w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng.InsertXML($xml)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
